$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("I6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("K6").Value = "7,97 TL - 15,96 TL - 199,41 TL"

# Row 12
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
$ws.Range("K12").Value = "WU: ,USD–; Diğer: 529 TL–4.454,74 TL"

# Row 13
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 795 TL | Azami 4.005 TL"
$ws.Range("I13").Value = "Hesaba: Asgari 1 TL | Azami 7,97 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14
$ws.Range("D14").Value = "3.500 TL - 13.500 TL"
$ws.Range("F14").Value = "2.785,72 TL - 12.380,95 TL"
$ws.Range("K14").Value = "1.196,51 TL - 5.583,74 TL"

# Row 24
$ws.Range("D24").Value = ""

# Row 25
$ws.Range("D25").Value = ""
